$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Feb 4 ("L" column) order amounts for the retailers that received
# orders that day. The sheet's F (Amount to collect) and G (Total Order,
# SUM(I:AM)) columns are formulas, so they recalculate automatically, as do
# the row-2 grand totals.
$l4Updates = @{
    7  = 2080
    10 = 5200
    19 = 2080
    23 = 3120
    24 = 5200
    31 = 2080
    32 = 1040
    40 = 2080
    41 = 3120
    46 = 3120
    48 = 3120
    52 = 1040
    56 = 5200
    60 = 2080
    62 = 2080
    71 = 3120
}

foreach ($row in $l4Updates.Keys) {
    $ws.Cells.Item($row, 12).Value = $l4Updates[$row]
}

# Row 53 got highlighted (tan/brown fill) by the author while entering it.
$ws.Range("L53").Value = 2080
$ws.Range("L53").Interior.Color = 8698081

# Restore the active selection to where the author left off.
$ws.Range("L51").Select()
